# Update cryptos list with latest scraped prices/volumes (GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "62.011.33"
$ws.Range("E2").Value = "  -1.87%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.422.13"
$ws.Range("E3").Value = "  -1.39%  "

# Row 4 - TetherUSD (leading apostrophe keeps price as literal text, not a number)
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.01%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'578.14"
$ws.Range("E5").Value = "  -0.55%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'152.75"
$ws.Range("E6").Value = "  +4.07%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  +1.37%  "

# Row 9 - Toncoin
$ws.Range("E9").Value = "  +4.68%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -0.51%  "

# Row 11 - Cardano
$ws.Range("E11").Value = "  +3.36%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "4.005.72"
$ws.Range("E12").Value = "  -1.48%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +0.74%  "

# Row 14 - Avalanche
$ws.Range("D14").Value = "'28.75"
$ws.Range("E14").Value = "  -2.32%  "

# Row 15 - WrappedEther
$ws.Range("D15").Value = "3.424.26"
$ws.Range("E15").Value = "  -1.52%  "

# Row 16 - ShibaInu
$ws.Range("E16").Value = "  -0.11%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "62.096.82"
$ws.Range("E17").Value = "  -1.82%  "

# Row 18 - Polkadot
$ws.Range("E18").Value = "  +1.60%  "

# Row 19
$ws.Range("E19").Value = "  +0.62%  "

# Row 20
$ws.Range("D20").Value = "'8.95"
$ws.Range("E20").Value = "  -4.07%  "

# Row 21
$ws.Range("D21").Value = "'383.25"
$ws.Range("E21").Value = "  -1.27%  "

# Row 22
$ws.Range("E22").Value = "  +0.83%  "

# Row 23
$ws.Range("D23").Value = "'75.29"
$ws.Range("E23").Value = "  +0.89%  "

# Row 24
$ws.Range("E24").Value = "  +0.03%  "

# Row 25
$ws.Range("D25").Value = "3.562.79"
$ws.Range("E25").Value = "  -1.61%  "

# Row 26
$ws.Range("E26").Value = "  -3.29%  "

# Row 27
$ws.Range("E27").Value = "  -0.66%  "

# Row 28
$ws.Range("E28").Value = "  +1.10%  "

# Row 29
$ws.Range("E29").Value = "  +0.01%  "

# Row 30 & 31 swap places: InternetComputer(DFINITY) <-> PancakeSwap
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "'2.12"
$ws.Range("E30").Value = "  -0.99%  "

$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "'7.90"
$ws.Range("E31").Value = "  -4.18%  "

# Row 32 - USDe
$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = "  -0.03%  "

# Row 33
$ws.Range("E33").Value = "  -2.20%  "

# Row 34
$ws.Range("E34").Value = "  -1.09%  "

# Row 35 - NEARProtocol
$ws.Range("D35").Value = "'5.53"
$ws.Range("E35").Value = "  +3.43%  "

# Row 36 - ImmutableX
$ws.Range("D36").Value = "'1.62"
$ws.Range("E36").Value = "  +0.71%  "

# Row 37
$ws.Range("E37").Value = "  -3.15%  "

# Row 38 - EnergySwap
$ws.Range("D38").Value = "'31.22"
$ws.Range("E38").Value = "  -2.35%  "

# Row 39 - Monero
$ws.Range("D39").Value = "'168.46"
$ws.Range("E39").Value = "  -0.32%  "

# Row 40
$ws.Range("E40").Value = "  -1.61%  "

# Row 41
$ws.Range("E41").Value = "  +2.64%  "

# Row 42 - OKB
$ws.Range("D42").Value = "'42.76"
$ws.Range("E42").Value = "  +1.03%  "

# Row 43
$ws.Range("E43").Value = "  -2.49%  "

# Row 44 - Filecoin
$ws.Range("D44").Value = "'4.41"
$ws.Range("E44").Value = "  +0.35%  "

# Row 45
$ws.Range("E45").Value = "  -3.15%  "

# Row 46
$ws.Range("E46").Value = "  -2.85%  "

# Row 47 - Maker (only price changes)
$ws.Range("D47").Value = "2.551.30"

# Row 48
$ws.Range("E48").Value = "  +0.69%  "

# Row 49
$ws.Range("E49").Value = "  -4.62%  "

# Row 50 - InjectiveProtocol
$ws.Range("D50").Value = "'22.60"
$ws.Range("E50").Value = "  -2.22%  "

# Row 51
$ws.Range("E51").Value = "  +0.13%  "
